# Master_Fuel_Sector_List.xlsx edit
# Commit message: "Remove older pulp and paper and smelting activity data"
#
# The "Sectors" sheet has rows for 2C_Metal-production and
# 2H_Pulp-and-paper-food-beverage-wood that previously carried bespoke
# "activity"/"units"/"description" values (Metal_S_Content / kt / Default
# total S content of smelted metals, and Pulp_Paper_Production / tons).
# Those bespoke values are removed so the rows fall back to the common
# GDP / B2005USD activity used by the rest of the sheet, and the now-unused
# shared strings disappear from the workbook on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# Row 37: 2C_Metal-production
$ws.Cells.Item(37, 2).Value = "GDP"
$ws.Cells.Item(37, 3).Value = "B2005USD"
$ws.Cells.Item(37, 4).Value = ""

# Row 42: 2H_Pulp-and-paper-food-beverage-wood
$ws.Cells.Item(42, 2).Value = "GDP"
$ws.Cells.Item(42, 3).Value = "B2005USD"

# Restore the view to where the user last left it (scrolled down, B30
# selected) rather than the original top-of-sheet selection.
$ws.Range("B30").Select()
